$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.506.11'
$ws.Range('E2').Value = '  +3.66%  '
$ws.Range('D3').Value = '1.587.00'
$ws.Range('E3').Value = '  +1.03%  '
$ws.Range('E4').Value = '  +0.97%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '212.98'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.94%  '
$ws.Range('E6').Value = '  +0.20%  '
$ws.Range('E7').Value = '  +1.00%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '24.22'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +5.52%  '
$ws.Range('E9').Value = '  +0.60%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0601'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.96%  '
$ws.Range('E11').Value = '  +1.63%  '
$ws.Range('D12').Value = '1.813.42'
$ws.Range('E12').Value = '  +1.04%  '
$ws.Range('D13').Value = '1.585.18'
$ws.Range('E13').Value = '  +0.88%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.529'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.88%  '
$ws.Range('D16').Value = '28.530.64'
$ws.Range('E16').Value = '  +3.89%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '63.08'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.18%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '230.83'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +2.28%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.50'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.13%  '
$ws.Range('E21').Value = '  +0.94%  '
$ws.Range('E22').Value = '  -1.73%  '
$ws.Range('E23').Value = '  -0.92%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.99'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.28%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '151.93'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.09%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '15.22'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.44%  '
$ws.Range('E27').Value = '  -0.78%  '
$ws.Range('E28').Value = '  -0.66%  '
$ws.Range('E29').Value = '  +0.94%  '
$ws.Range('E30').Value = '  -0.85%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0470'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.66%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.26'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.28%  '
$ws.Range('E33').Value = '  +1.51%  '
$ws.Range('D34').Value = '1.395.11'
$ws.Range('E34').Value = '  -4.09%  '
$ws.Range('E35').Value = '  -1.42%  '
$ws.Range('E36').Value = '  -10.53%  '
$ws.Range('E37').Value = '  +1.04%  '
$ws.Range('E38').Value = '  +10.65%  '
$ws.Range('E39').Value = '  -1.05%  '
$ws.Range('E40').Value = '  -0.17%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.812'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.11%  '
$ws.Range('E42').Value = '  +0.92%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.63'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.40%  '
$ws.Range('E44').Value = '  +0.85%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.980'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.89%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '62.95'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.33%  '
$ws.Range('D47').Value = '1.723.17'
$ws.Range('E47').Value = '  +0.93%  '
$ws.Range('E48').Value = '  +0.95%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '86.61'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.25%  '
$ws.Range('E50').Value = '  -0.77%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0520'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.95%  '
